{"js": "// Update the answer table: replace each filled cell's text with the new\n// division-problem answer, matched strictly by (row, column) position so\n// that duplicate text values (e.g. \"57\u00f74=14, 1\" appearing twice) are not\n// ambiguous.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of row index -> array of 5 new cell texts (only rows that contain\n// data are listed; the blank spacer rows are left untouched).\nconst updates = {\n  0: [\"77\u00f77=11, 0\", \"34\u00f76=5, 4\", \"60\u00f79=6, 6\", \"14\u00f74=3, 2\", \"76\u00f73=25, 1\"],\n  4: [\"49\u00f72=24, 1\", \"86\u00f77=12, 2\", \"23\u00f77=3, 2\", \"69\u00f76=11, 3\", \"96\u00f72=48, 0\"],\n  8: [\"73\u00f79=8, 1\", \"67\u00f74=16, 3\", \"54\u00f72=27, 0\", \"19\u00f75=3, 4\", \"93\u00f76=15, 3\"],\n  12: [\"26\u00f75=5, 1\", \"24\u00f75=4, 4\", \"42\u00f75=8, 2\", \"13\u00f74=3, 1\", \"69\u00f77=9, 6\"],\n  16: [\"84\u00f77=12, 0\", \"94\u00f75=18, 4\", \"14\u00f76=2, 2\", \"60\u00f78=7, 4\", \"84\u00f73=28, 0\"],\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const newTexts = updates[rowIndexStr];\n  for (let colIndex = 0; colIndex < newTexts.length; colIndex++) {\n    const cell = table.getCellOrNullObject(rowIndex, colIndex);\n    cell.value = newTexts[colIndex];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the answer table: replace each filled cell's text with the new\n# division-problem answer, matched strictly by (row, column) position so\n# that duplicate text values (e.g. \"57\u00f74=14, 1\" appearing twice) are not\n# ambiguous.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Hashtable of 1-based row index -> array of 5 new cell texts (only rows\n# that contain data are listed; the blank spacer rows are left untouched).\n$updates = @{\n    1  = @(\"77\u00f77=11, 0\", \"34\u00f76=5, 4\", \"60\u00f79=6, 6\", \"14\u00f74=3, 2\", \"76\u00f73=25, 1\")\n    5  = @(\"49\u00f72=24, 1\", \"86\u00f77=12, 2\", \"23\u00f77=3, 2\", \"69\u00f76=11, 3\", \"96\u00f72=48, 0\")\n    9  = @(\"73\u00f79=8, 1\", \"67\u00f74=16, 3\", \"54\u00f72=27, 0\", \"19\u00f75=3, 4\", \"93\u00f76=15, 3\")\n    13 = @(\"26\u00f75=5, 1\", \"24\u00f75=4, 4\", \"42\u00f75=8, 2\", \"13\u00f74=3, 1\", \"69\u00f77=9, 6\")\n    17 = @(\"84\u00f77=12, 0\", \"94\u00f75=18, 4\", \"14\u00f76=2, 2\", \"60\u00f78=7, 4\", \"84\u00f73=28, 0\")\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $newTexts = $updates[$rowIndex]\n    for ($colIndex = 1; $colIndex -le $newTexts.Length; $colIndex++) {\n        $cell = $t.Cell($rowIndex, $colIndex)\n        $cell.Range.Text = $newTexts[$colIndex - 1]\n    }\n}\n"}
